# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) for the leve rows whose market-board prices changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 681.4
$ws.Range("I38").Value = 296.25
$ws.Range("J38").Value = 2222
$ws.Range("K38").Value = 888.75
$ws.Range("L38").Value = 6666
$ws.Range("M38").Value = -516.75
$ws.Range("N38").Value = -7410
$ws.Range("H39").Value = 999.5
$ws.Range("I39").Value = 399.4
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 1198.2
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = -902.1999999999998
$ws.Range("N39").Value = -12592
$ws.Range("H40").Value = 2887.9092
$ws.Range("I40").Value = 2755.4285
$ws.Range("J40").Value = 3119.75
$ws.Range("K40").Value = 2755.4285
$ws.Range("L40").Value = 3119.75
$ws.Range("M40").Value = -2580.4285
$ws.Range("N40").Value = -3469.75
$ws.Range("H43").Value = 4036.4
$ws.Range("J43").Value = 4249.6924
$ws.Range("L43").Value = 4249.6924
$ws.Range("N43").Value = -4387.6924
$ws.Range("H88").Value = 2728.5833
$ws.Range("I88").Value = 1935.75
$ws.Range("J88").Value = 3125
$ws.Range("K88").Value = 1935.75
$ws.Range("L88").Value = 3125
$ws.Range("M88").Value = -1529.75
$ws.Range("N88").Value = -3937
$ws.Range("H91").Value = 2728.5833
$ws.Range("I91").Value = 1935.75
$ws.Range("J91").Value = 3125
$ws.Range("K91").Value = 1935.75
$ws.Range("L91").Value = 3125
$ws.Range("M91").Value = -531.75
$ws.Range("N91").Value = -5933
$ws.Range("H113").Value = 15543.667
$ws.Range("I113").Value = 4649.5
$ws.Range("J113").Value = 20990.75
$ws.Range("K113").Value = 4649.5
$ws.Range("L113").Value = 20990.75
$ws.Range("M113").Value = -1395.5
$ws.Range("N113").Value = -27498.75
$ws.Range("H129").Value = 19616854
$ws.Range("I129").Value = 3198.1428
$ws.Range("J129").Value = 33346414
$ws.Range("K129").Value = 9594.428400000001
$ws.Range("L129").Value = 100039242
$ws.Range("M129").Value = -4594.428400000001
$ws.Range("N129").Value = -100049242
$ws.Range("H132").Value = 1178628.9
$ws.Range("I132").Value = 1836774
$ws.Range("K132").Value = 5510322
$ws.Range("M132").Value = -5507792
$ws.Range("H135").Value = 4748.077
$ws.Range("J135").Value = 6571.143
$ws.Range("L135").Value = 59140.287
$ws.Range("N135").Value = -64210.287
$ws.Range("H137").Value = 3411.3044
$ws.Range("I137").Value = 2618.8462
$ws.Range("K137").Value = 7856.5386
$ws.Range("M137").Value = -5306.5386
$ws.Range("H141").Value = 3050.8
$ws.Range("I141").Value = 3050.8
$ws.Range("K141").Value = 9152.400000000001
$ws.Range("M141").Value = -3972.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5527.84
$ws.Range("I2").Value = 5540.9414
$ws.Range("J2").Value = 5500
$ws.Range("K2").Value = 5540.9414
$ws.Range("L2").Value = 5500
$ws.Range("M2").Value = -5427.9414
$ws.Range("N2").Value = -5726
$ws.Range("H116").Value = 5527.84
$ws.Range("I116").Value = 5540.9414
$ws.Range("J116").Value = 5500
$ws.Range("K116").Value = 5540.9414
$ws.Range("L116").Value = 5500
$ws.Range("M116").Value = -3246.9414
$ws.Range("N116").Value = -10088
$ws.Range("H132").Value = 3901429
$ws.Range("J132").Value = 180998
$ws.Range("L132").Value = 542994
$ws.Range("N132").Value = -548054

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5527.84
$ws.Range("I3").Value = 5540.9414
$ws.Range("J3").Value = 5500
$ws.Range("K3").Value = 5540.9414
$ws.Range("L3").Value = 5500
$ws.Range("M3").Value = -5426.9414
$ws.Range("N3").Value = -5728
$ws.Range("H86").Value = 5425.0713
$ws.Range("J86").Value = 6682.857
$ws.Range("L86").Value = 6682.857
$ws.Range("N86").Value = -8928.857
$ws.Range("H89").Value = 5425.0713
$ws.Range("J89").Value = 6682.857
$ws.Range("L89").Value = 33414.285
$ws.Range("N89").Value = -44646.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 38467330
$ws.Range("I16").Value = 83336540
$ws.Range("K16").Value = 83336540
$ws.Range("M16").Value = -83336253
$ws.Range("H58").Value = 83344130
$ws.Range("I58").Value = 166674770
$ws.Range("K58").Value = 166674770
$ws.Range("M58").Value = -166674567
$ws.Range("H113").Value = 38467330
$ws.Range("I113").Value = 83336540
$ws.Range("K113").Value = 83336540
$ws.Range("M113").Value = -83334370
$ws.Range("H132").Value = 6211.1924
$ws.Range("I132").Value = 5054.55
$ws.Range("J132").Value = 10066.667
$ws.Range("K132").Value = 15163.65
$ws.Range("L132").Value = 30200.001
$ws.Range("M132").Value = -12633.65
$ws.Range("N132").Value = -35260.001
$ws.Range("H134").Value = 33342358
$ws.Range("J134").Value = 11778.048
$ws.Range("L134").Value = 35334.144
$ws.Range("N134").Value = -40404.144
$ws.Range("H136").Value = 83344130
$ws.Range("I136").Value = 166674770
$ws.Range("K136").Value = 500024310
$ws.Range("M136").Value = -500021760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 109879.586
$ws.Range("I121").Value = 286.25
$ws.Range("K121").Value = 858.75
$ws.Range("M121").Value = 451.25
$ws.Range("H122").Value = 142348.78
$ws.Range("J122").Value = 155228.84
$ws.Range("L122").Value = 1397059.56
$ws.Range("N122").Value = -1401959.56

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 8774.467000000001
$ws.Range("J113").Value = 13223.444
$ws.Range("L113").Value = 13223.444
$ws.Range("N113").Value = -17563.444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 673.2857
$ws.Range("I22").Value = 652.8
$ws.Range("K22").Value = 652.8
$ws.Range("M22").Value = -357.8
$ws.Range("H27").Value = 673.2857
$ws.Range("I27").Value = 652.8
$ws.Range("K27").Value = 652.8
$ws.Range("M27").Value = -545.8
$ws.Range("H46").Value = 23810378
$ws.Range("I46").Value = 819.9375
$ws.Range("J46").Value = 100000960
$ws.Range("K46").Value = 819.9375
$ws.Range("L46").Value = 100000960
$ws.Range("M46").Value = -631.9375
$ws.Range("N46").Value = -100001336
$ws.Range("H55").Value = 2259.724
$ws.Range("I55").Value = 1122.4286
$ws.Range("J55").Value = 3321.2
$ws.Range("K55").Value = 1122.4286
$ws.Range("L55").Value = 3321.2
$ws.Range("M55").Value = -949.4286
$ws.Range("N55").Value = -3667.2
$ws.Range("H132").Value = 3271.25
$ws.Range("I132").Value = 3125.0667
$ws.Range("J132").Value = 3709.8
$ws.Range("K132").Value = 9375.2001
$ws.Range("L132").Value = 11129.4
$ws.Range("M132").Value = -6845.2001
$ws.Range("N132").Value = -16189.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 785.5
$ws.Range("I81").Value = 749.375
$ws.Range("K81").Value = 1498.75
$ws.Range("M81").Value = -437.75
$ws.Range("H84").Value = 785.5
$ws.Range("I84").Value = 749.375
$ws.Range("K84").Value = 7493.75
$ws.Range("M84").Value = -2189.75
$ws.Range("H100").Value = 1758.25
$ws.Range("I100").Value = 1043
$ws.Range("K100").Value = 2086
$ws.Range("M100").Value = -1545
